$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "ling"
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "title"
$ws.Range("D1").Value = "department"
$ws.Range("E1").Value = "location"
$ws.Range("F1").Value = "deadline"
$ws.Range("G1").Value = "number"
$ws.Range("H1").Value = "post date"
$ws.Range("I1").Value = "interlinkregex"
$ws.Range("J1").Value = "finallinkregex"

# --- Row 2 (cmbchina, already present - refresh values / hyperlink stays) ---
$ws.Range("A2").Value = "http://career.cmbchina.com/Campus/Campus.aspx"
$ws.Range("B2").Value = "cmbchina"
$ws.Range("C2").Value = "//*[@id=`"rightdiv`"]/div/div[1]/text()"
$ws.Range("D2").Value = "//*[@id=`"rightdiv`"]/div/div[2]/p[3]/text()"
$ws.Range("E2").Value = "//*[@id=`"rightdiv`"]/div/div[2]/p[5]/text()"
$ws.Range("F2").Value = "//*[@id=`"rightdiv`"]/div/div[2]/p[6]/text()"
$ws.Range("I2").Value = "branch="
$ws.Range("J2").Value = "Position.aspx.id"

# --- Row 3 (ccb - new) ---
$ws.Range("A3").Value = "http://job.ccb.com/ccbjob/cn/job/index.jsp"
$ws.Hyperlinks.Add($ws.Range("A3"), "http://job.ccb.com/ccbjob/cn/job/index.jsp")
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("B3").Value = "ccb"
$ws.Range("C3").Value = "//*[@id=`"title`"]/strong/text()"
$ws.Range("D3").Value = "//*[@id=`"title`"]/strong/text()"
$ws.Range("E3").Value = "//*[@id=`"data`"]/table/tbody/tr[2]/td/p[5]/font/text()"
$ws.Range("H3").Value = "//*[@id=`"data`"]/table[1]/tbody/tr[1]/td/text()"
$ws.Range("I3").Value = "branch_notice_list2"
$ws.Range("J3").Value = "info"

# --- Row 4 (abchina - new) ---
$ws.Range("A4").Value = "http://job.abchina.com/rio/index.do?action=openHome"
$ws.Hyperlinks.Add($ws.Range("A4"), "http://job.abchina.com/rio/index.do?action=openHome")
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("B4").Value = "abchina"
$ws.Range("C4").Value = "//*[@id=`"frmPutQuestion`"]/table[3]/tr[2]/td[2]"
$ws.Range("D4").Value = "//*[@id=`"frmPutQuestion`"]/table[3]/tr[6]/td[2]"
$ws.Range("E4").Value = "//*[@id=`"frmPutQuestion`"]/table[3]/tbody/tr[3]/td[2]"
$ws.Range("F4").Value = "//*[@id=`"frmPutQuestion`"]/table[3]/tr[5]/td[2]"
$ws.Range("G4").Value = "//*[@id=`"frmPutQuestion`"]/table[3]/tr[2]/td[4]"
$ws.Range("H4").Value = "//*[@id=`"frmPutQuestion`"]/table[3]/tr[1]/td[4]"
$ws.Range("I4").Value = "openHome"
$ws.Range("J4").Value = "jobDetails"

# --- Row 6 (sample array, unchanged content but shifts shared-string index) ---
$ws.Range("B6").Value = "['http://career.cmbchina.com/Campus/Position.aspx?id=10234', 'cmbchina', 'IT岗', '无锡分行', '无锡分行及下辖机构', '2016-10-14', 'Not available', 'Not available']"

# --- Column widths ---
$ws.Columns.Item(3).ColumnWidth = 15.5
$ws.Columns.Item(6).ColumnWidth = 35.33

# --- Selection ---
$ws.Range("F12").Select()
